# This script reproduces the cryptocurrency price/volume refresh described
# by the commit: dozens of individual cell updates in the "Price" (D) and
# "Volume(1h)" (E) columns of the single data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A handful of the new Price values (e.g. "86.93") are plain decimals that
# Excel would otherwise auto-convert to a Number when assigned through
# Range.Value. The source data stores every Price/Volume cell as literal
# text, so force those specific cells to Text format first (looping over
# single-cell ranges -- a multi-area Range("A1,A2,...") only honours the
# first area for NumberFormat/Style in this host). Values with extra
# punctuation (thousands separators, subscripts, etc.) already fail numeric
# parsing on their own and do not need this treatment.
$textForceRefs = @('D5', 'D6', 'D9', 'D10', 'D11', 'D12', 'D18', 'D21', 'D22', 'D23', 'D24', 'D26', 'D27', 'D28', 'D30', 'D31', 'D32', 'D33', 'D34', 'D35', 'D37', 'D40', 'D42', 'D46', 'D47', 'D48', 'D50')
foreach ($ref in $textForceRefs) {
  $ws.Range($ref).NumberFormat = "@"
}

$ws.Range('D2').Value = '40.055.99'
$ws.Range('E2').Value = '  +2.59%  '
$ws.Range('D3').Value = '2.242.91'
$ws.Range('E3').Value = '  -0.12%  '
$ws.Range('E4').Value = '  -0.24%  '
$ws.Range('D5').Value = '293.98'
$ws.Range('E5').Value = '  +0.34%  '
$ws.Range('D6').Value = '86.93'
$ws.Range('E6').Value = '  +8.78%  '
$ws.Range('E7').Value = '  +2.21%  '
$ws.Range('E8').Value = '  -0.16%  '
$ws.Range('D9').Value = '0.475'
$ws.Range('E9').Value = '  +4.37%  '
$ws.Range('D10').Value = '30.95'
$ws.Range('E10').Value = '  +11.92%  '
$ws.Range('D11').Value = '0.0804'
$ws.Range('E11').Value = '  +4.97%  '
$ws.Range('D12').Value = '46.97'
$ws.Range('E12').Value = '  +5.23%  '
$ws.Range('E13').Value = '  +0.11%  '
$ws.Range('E14').Value = '  +7.46%  '
$ws.Range('D15').Value = '2.584.89'
$ws.Range('E15').Value = '  -0.43%  '
$ws.Range('E16').Value = '  +2.47%  '
$ws.Range('D17').Value = '2.242.12'
$ws.Range('E17').Value = '  -1.02%  '
$ws.Range('D18').Value = '0.734'
$ws.Range('E18').Value = '  +3.58%  '
$ws.Range('D19').Value = '39.970.46'
$ws.Range('E19').Value = '  +2.57%  '
$ws.Range('D20').Value = '0.0₃0898'
$ws.Range('E20').Value = '  +5.34%  '
$ws.Range('D21').Value = '5.85'
$ws.Range('E21').Value = '  +2.12%  '
$ws.Range('D22').Value = '10.64'
$ws.Range('E22').Value = '  +8.37%  '
$ws.Range('D23').Value = '65.70'
$ws.Range('E23').Value = '  +1.42%  '
$ws.Range('D24').Value = '236.49'
$ws.Range('E24').Value = '  +5.05%  '
$ws.Range('E25').Value = '  +0.08%  '
$ws.Range('D26').Value = '2.45'
$ws.Range('E26').Value = '  +3.44%  '
$ws.Range('D27').Value = '1.84'
$ws.Range('E27').Value = '  +6.79%  '
$ws.Range('D28').Value = '23.11'
$ws.Range('E28').Value = '  +4.71%  '
$ws.Range('E29').Value = '  +2.12%  '
$ws.Range('D30').Value = '9.29'
$ws.Range('E30').Value = '  +6.11%  '
$ws.Range('D31').Value = '34.25'
$ws.Range('E31').Value = '  +10.34%  '
$ws.Range('D32').Value = '154.72'
$ws.Range('E32').Value = '  +4.78%  '
$ws.Range('D33').Value = '0.999'
$ws.Range('E33').Value = '  -0.29%  '
$ws.Range('D34').Value = '4.89'
$ws.Range('E34').Value = '  +4.03%  '
$ws.Range('D35').Value = '0.0715'
$ws.Range('E35').Value = '  +5.62%  '
$ws.Range('E36').Value = '  +0.71%  '
$ws.Range('D37').Value = '16.64'
$ws.Range('E37').Value = '  +14.57%  '
$ws.Range('E38').Value = '  +6.87%  '
$ws.Range('E39').Value = '  +2.92%  '
$ws.Range('D40').Value = '2.72'
$ws.Range('E40').Value = '  +4.61%  '
$ws.Range('E41').Value = '  +6.55%  '
$ws.Range('D42').Value = '3.82'
$ws.Range('E42').Value = '  +6.53%  '
$ws.Range('D43').Value = '1.968.91'
$ws.Range('E43').Value = '  +2.49%  '
$ws.Range('E44').Value = '  -0.28%  '
$ws.Range('E45').Value = '  +7.66%  '
$ws.Range('D46').Value = '9.71'
$ws.Range('E46').Value = '  +8.27%  '
$ws.Range('D47').Value = '16.36'
$ws.Range('E47').Value = '  +1.75%  '
$ws.Range('D48').Value = '2.60'
$ws.Range('E48').Value = '  +4.04%  '
$ws.Range('D49').Value = '2.457.94'
$ws.Range('E49').Value = '  -1.26%  '
$ws.Range('D50').Value = '71.25'
$ws.Range('E50').Value = '  +8.38%  '
$ws.Range('E51').Value = '  +16.15%  '

# Restore the default "Normal" style on the cells we temporarily forced to
# Text above, so only the values (not the formatting) differ from the
# original workbook.
foreach ($ref in $textForceRefs) {
  $ws.Range($ref).Style = "Normal"
}

